# edit.ps1
#
# Reshapes the "Output" report so that two new columns appear:
#   - "Numar zile"    inserted right after "Perioada Internarii" (pushes the
#                      old "Urgenta" column one slot to the right, F -> G)
#   - "Hipertensiune"  inserted right after "Urgenta" (pushes the old
#                      "LDL COLESTEROL" column one slot to the right, G -> I)
#
# Original layout: A Nume | B Prenume | C Varsta | D Data Tiparire |
#                   E Perioada Internarii | F Urgenta | G LDL COLESTEROL
#
# New layout:       A Nume | B Prenume | C Varsta | D Data Tiparire |
#                   E Perioada Internarii | F Numar zile | G Urgenta |
#                   H Hipertensiune | I LDL COLESTEROL
#
# The "Perioada Internarii" text is also trimmed down to just the start/end
# dates (the embedded end-time and day-count are dropped - the day count now
# lives in its own "Numar zile" column), and the new "Numar zile" /
# "Hipertensiune" / "LDL COLESTEROL" values are filled in for the five
# patient rows. Columns A-D are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) -------------------------------------------------------
# I1 and H1 are brand-new header cells - copy the header style (bold, boxed,
# centered) from an existing header cell before writing their text so they
# pick up the same cell-format index as the rest of row 1.
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(1, 8).PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(1, 9).Value = "LDL COLESTEROL"   # I1 - was G1
$ws.Cells.Item(1, 8).Value = "Hipertensiune"    # H1 - new
$ws.Cells.Item(1, 7).Value = "Urgenta"          # G1 - was F1
$ws.Cells.Item(1, 6).Value = "Numar zile"       # F1 - new

# --- Patient data rows (2-6) ------------------------------------------------
# Row -> Perioada Internarii (trimmed), Numar zile, Urgenta, Hipertensiune, LDL COLESTEROL
$rows = @(
    @{ Row = 2; Period = "26/11/2019 08:04 - 02/12/2019"; Zile = "6"; Urgenta = "NU `n"; Hiper = "Da"; Ldl = "97.36000000000001" },
    @{ Row = 3; Period = "05/12/2019 17:03 - 13/12/2019"; Zile = "8"; Urgenta = "DA `n"; Hiper = "Da"; Ldl = "49.72000000000001" },
    @{ Row = 4; Period = "18/11/2019 09:20 - 22/11/2019"; Zile = "4"; Urgenta = "NU `n"; Hiper = "Da"; Ldl = "97.44" },
    @{ Row = 5; Period = "25/11/2019 09:19 - 29/11/2019"; Zile = "4"; Urgenta = "NU `n"; Hiper = "Da"; Ldl = "114.38" },
    @{ Row = 6; Period = "19/11/2019 08:23 - 22/11/2019"; Zile = "3"; Urgenta = "NU `n"; Hiper = "Da"; Ldl = "93.82" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 5).Value = $r.Period     # E - Perioada Internarii (trimmed)

    # F (Numar zile) and I (LDL COLESTEROL) hold digit strings ("6",
    # "97.36000000000001", ...) that must stay text, not auto-convert to a
    # number. Format as text before writing so the value is stored as a
    # string, then reset the style back to Normal so no visible formatting
    # (or extra style) lingers on the cell.
    $ws.Cells.Item($n, 6).NumberFormat = "@"
    $ws.Cells.Item($n, 6).Value = $r.Zile       # F - Numar zile (new)
    $ws.Cells.Item($n, 6).Style = "Normal"

    $ws.Cells.Item($n, 7).Value = $r.Urgenta    # G - Urgenta (was F)
    $ws.Cells.Item($n, 8).Value = $r.Hiper      # H - Hipertensiune (new)

    $ws.Cells.Item($n, 9).NumberFormat = "@"
    $ws.Cells.Item($n, 9).Value = $r.Ldl        # I - LDL COLESTEROL (new)
    $ws.Cells.Item($n, 9).Style = "Normal"
}

# --- Trailing empty template rows (7-51) ------------------------------------
# These rows had no data in F/G (Urgenta/LDL) before, and keep it that way -
# the only change is that the grid now extends two columns further right
# (through I) to match the new header row. A blank cell has no Value to
# write, so the empty placeholder cells are materialised by copying the
# (blank) format of column A down into the new H/I columns instead.
$ws.Range("A7:A51").Copy()
$ws.Range("H7:H51").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A7:A51").Copy()
$ws.Range("I7:I51").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Undo the automatic row-height bump that Excel applies when a multi-line
# value (embedded "\n") is written into a cell, so rows keep their original
# (implicit/default) height instead of gaining an explicit customHeight.
$ws.Range("A1:A51").EntireRow.AutoFit()
